# CompStat weekly update: Volume/Number, report date range, and the
# Week-to-Date / 28-Day / Year-to-Date / 2-Year crime figures (rows 15-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 30   Number  15" -> "...16" and the report date
# range "4/10/2023 ... 4/16/2023" -> "4/17/2023 ... 4/23/2023".
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# ---------------------------------------------------------------------
# Helper pattern used below for cells that flip between a numeric value
# and the literal text placeholders "0" / "***.*" (shared strings used
# for "no data" cells). Plain `.Value = "0"` would be auto-coerced back
# to a number by Excel, so we stage the value as Text first and then
# paste-special the number format back from a donor cell that already
# carries the desired style (General / #,##0 / #,##0.0) - this keeps
# the cell's font/alignment/style intact while only changing the
# number format + content type.
# ---------------------------------------------------------------------

# Donor cells with the three styles we need to re-apply:
#   C14 -> General (text placeholders "0" / "***.*")
#   I14 -> #,##0   (plain integer counts)
#   K14 -> #,##0.0 (percent-change figures)

# Row 15 - Rape: F15 numeric 1 -> text "0"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)

# Row 16 - Robbery
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -64.705882352941
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -15.384615384615
$ws.Range("M16").Value = -31.25

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -4.545454545454
$ws.Range("I17").Value = 81
$ws.Range("J17").Value = 89
$ws.Range("K17").Value = -8.988764044943
$ws.Range("L17").Value = -8.988764044943
$ws.Range("M17").Value = 55.769230769230

# Row 18 - Burglary: C18 text "0" -> numeric 1
$ws.Range("I14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -75
$ws.Range("I18").Value = 46
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = 31.428571428571
$ws.Range("L18").Value = -24.590163934426
$ws.Range("M18").Value = 53.333333333333

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 25
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 92
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = -8
$ws.Range("L19").Value = -3.157894736842
$ws.Range("M19").Value = 19.480519480519

# Row 20 - G.L.A.: C20 text "0" -> numeric 3, D20 numeric 3 -> text "0",
# E20 numeric -100 -> text "***.*"
$ws.Range("I14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 3

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 40
$ws.Range("K20").Value = -9.090909090909
$ws.Range("L20").Value = 166.666666666667
$ws.Range("M20").Value = 73.913043478260

# Row 21 - TOTAL
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -11.111111111111
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -30.337078651685
$ws.Range("I21").Value = 307
$ws.Range("J21").Value = 327
$ws.Range("K21").Value = -6.116207951070
$ws.Range("L21").Value = -2.229299363057
$ws.Range("M21").Value = 21.825396825396

# Row 22 - Transit: C22 numeric 1 -> text "0"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = -75
$ws.Range("L22").Value = -25

# Row 23 - Housing
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = -76.923076923076
$ws.Range("L23").Value = -50
$ws.Range("M23").Value = -40

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -41.379310344827
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 138
$ws.Range("H24").Value = -37.681159420289
$ws.Range("I24").Value = 300
$ws.Range("J24").Value = 509
$ws.Range("K24").Value = -41.060903732809
$ws.Range("L24").Value = 3.806228373702
$ws.Range("M24").Value = 104.081632653061

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -63.636363636363
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = -24.137931034482
$ws.Range("I25").Value = 123
$ws.Range("J25").Value = 116
$ws.Range("K25").Value = 6.034482758620
$ws.Range("L25").Value = 16.037735849056
$ws.Range("M25").Value = -3.149606299212

# Row 26 - UCR Rape*
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0

# Row 27 - Other Sex Crimes: D27 text "0" -> numeric 2,
# E27 text "***.*" -> numeric -50
$ws.Range("C27").Value = 1

$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2

$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -50

$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 13
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = -13.333333333333
$ws.Range("L27").Value = 62.5
